$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header fields
# ------------------------------------------------------------------
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must remain stored as TEXT (not be
# coerced into a number) exactly like in the original file. Typing the
# value straight into the cell makes Excel auto-detect it as a number,
# so we stage it in a scratch cell first, force it to Text there, copy
# the resulting (already-text) value into B3, and finally restore B3's
# original look (style) by pasting formats from B2 (same base style).
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "2570314725427075"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z1").Clear()

$ws.Range("C3").Value = "Mohaupt"

# ------------------------------------------------------------------
# Opening balance line
# ------------------------------------------------------------------
$ws.Range("D5").Value = "KONTOSTAND AM 18.01.2025"

# ------------------------------------------------------------------
# Row 6 (transaction 1)
# ------------------------------------------------------------------
$ws.Range("B6").Value = "19.01."
$ws.Range("C6").Value = "20.01."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 6596850"
$ws.Range("E6").Value = "40,97-"

# ------------------------------------------------------------------
# Row 7 (transaction 2)
# ------------------------------------------------------------------
$ws.Range("B7").Value = "21.01."
$ws.Range("C7").Value = "22.01."
$ws.Range("D7").Value = "MCDONALDS Gießen"
$ws.Range("E7").Value = "24,25-"

# ------------------------------------------------------------------
# Row 8 (transaction 3)
# ------------------------------------------------------------------
$ws.Range("B8").Value = "22.01."
$ws.Range("C8").Value = "23.01."
$ws.Range("D8").Value = "PAYPAL PWWSAS"
$ws.Range("E8").Value = "70,05-"

# ------------------------------------------------------------------
# Rows 9-11 no longer carry a transaction: clear their contents.
# Column E's alignment differs between row 9 (centered) and rows
# 10-11 (right aligned); both combinations already exist elsewhere in
# the workbook's style table, so we build each one once in a scratch
# cell (seeded from an existing B-column cell that shares the same
# base font) and paste the resulting formats onto the target cells.
# This reproduces the exact look without inventing brand new styles.
# ------------------------------------------------------------------
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""

$ws.Range("Z1").Value = "x"
$ws.Range("B9").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("E9").Value = ""
$ws.Range("E9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z1").Clear()

$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""

$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""

$ws.Range("Z1").Value = "x"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z1").HorizontalAlignment = -4152   # xlRight
$ws.Range("Z1").VerticalAlignment = -4108     # xlCenter
$ws.Range("Z1").WrapText = $true
$ws.Range("Z1").Copy() | Out-Null

$ws.Range("E10").Value = ""
$ws.Range("E10").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E11").Value = ""
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("Z1").Clear()

# ------------------------------------------------------------------
# Closing balance line
# ------------------------------------------------------------------
$ws.Range("D12").Value = "KONTOSTAND AM 25.01.2025"
$ws.Range("E12").Value = "135,27-"

# ------------------------------------------------------------------
# Next billing date
# ------------------------------------------------------------------
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 31.01.2025"
